$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change B4 from "jYsebuv" to "Dhana@01" and turn it into a hyperlink
$ws.Range("B4").Value = "Dhana@01"
$ws.Hyperlinks.Add($ws.Range("B4"), "Dhana@01")

# Re-select a cell to mimic the post-edit cursor position recorded in the diff
$ws.Range("B7").Select()
